# Add basic CLI support / Fix error tests / Fix yml save
#
# excelcy test fixture update:
#  - "phase" sheet: the trailing save_nlp/export_entity rows are replaced by
#    two export_train calls (xlsx export + yml export).
#  - "train" sheet: the "span" column (offsets like "0,4") is dropped; the
#    "entity" column shifts left to take its place.
#  - "config" sheet: train_iteration goes from 2 to 10.

$wb = $excel.ActiveWorkbook

# --- "phase" sheet ---------------------------------------------------
$wsPhase = $wb.Worksheets.Item("phase")

# Row 5 used to be a bare save_nlp() call with no args; it becomes the
# first export_train() call, writing the xlsx.
$wsPhase.Range("C5").Value = "export_train"
$wsPhase.Range("D5").Value = "file_path=export/train_04.xlsx"

# Row 6 used to be export_entity(file_path=person.xlsx, label=PERSON); it
# becomes the second export_train() call, writing the yml.
$wsPhase.Range("C6").Value = "export_train"
$wsPhase.Range("D6").Value = "file_path=export/train_04.yml"

# --- "train" sheet -----------------------------------------------------
# Column D ("span", e.g. "0,4") is removed entirely; the old column E
# ("entity") shifts left into D. Deleting the whole column reproduces
# both the cell shift and the new A1:D16 dimension.
$wsTrain = $wb.Worksheets.Item("train")
$wsTrain.Columns.Item(4).Delete()

# --- "config" sheet -----------------------------------------------------
$wsConfig = $wb.Worksheets.Item("config")
$wsConfig.Range("B5").Value = 10

# Move the visible selection to B6 on the config sheet (matches the
# updated saved cursor position), then return focus to "phase" so it
# stays the active tab, as in the source workbook.
$wsConfig.Activate() | Out-Null
$wsConfig.Range("B6").Select() | Out-Null
$wsPhase.Activate() | Out-Null
